# Insert a new weekly price record as row 36 in the "Bruselas (repollito)"
# dataset. This shifts the existing rows 36-45 down to rows 37-46 and
# extends the used range from A1:R45 to A1:R46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36, pushing rows 36:45 down to 37:46.
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the latest weekly observation.
$ws.Cells.Item(36, 1).Value = 9
$ws.Cells.Item(36, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(36, 3).Value = "Metropolitana"
$ws.Cells.Item(36, 4).Value = 44754
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(36, 6).Value = 100112035
$ws.Cells.Item(36, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 43
$ws.Cells.Item(36, 11).Value = 22000
$ws.Cells.Item(36, 12).Value = 22000
$ws.Cells.Item(36, 13).Value = 22000
$ws.Cells.Item(36, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(36, 15).Value = "Hijuelas"
$ws.Cells.Item(36, 16).Value = 1467
$ws.Cells.Item(36, 17).Value = 15
$ws.Cells.Item(36, 18).Value = "Hortaliza"
